# Refresh the cryptocurrency price/volume snapshot (scheduled GitHub Actions scrape).
# Hedera and TrustWalletToken also swap ranking positions (rows 40-41).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D='22.396.12'; E='  +0.02%  ' }
    @{ Row=3; D='1.573.37'; E='  +0.01%  ' }
    @{ Row=4; E='  +0.05%  ' }
    @{ Row=5; E='  +0.13%  ' }
    @{ Row=6; D='291.06'; E='  +0.00%  ' }
    @{ Row=7; E='  +2.53%  ' }
    @{ Row=8; D='49.90'; E='  +0.58%  ' }
    @{ Row=9; D='0.3415'; E='  +1.38%  ' }
    @{ Row=10; D='1.163'; E='  -0.29%  ' }
    @{ Row=11; D='0.07665'; E='  +1.14%  ' }
    @{ Row=12; E='  +0.08%  ' }
    @{ Row=13; D='21.29'; E='  +0.81%  ' }
    @{ Row=14; D='5.975'; E='  -1.33%  ' }
    @{ Row=15; D='6.925'; E='  +0.86%  ' }
    @{ Row=16; D='1.575.95'; E='  -0.60%  ' }
    @{ Row=17; E='  -0.05%  ' }
    @{ Row=18; D='90.60'; E='  +1.28%  ' }
    @{ Row=19; D='0.06756'; E='  -0.08%  ' }
    @{ Row=20; E='  +0.18%  ' }
    @{ Row=21; E='  +2.65%  ' }
    @{ Row=22; D='6.233'; E='  +0.01%  ' }
    @{ Row=23; D='0.5282'; E='  -4.51%  ' }
    @{ Row=24; D='12.02'; E='  +0.33%  ' }
    @{ Row=25; D='22.404.50'; E='  -0.03%  ' }
    @{ Row=26; D='2.423'; E='  -0.12%  ' }
    @{ Row=27; D='2.759'; E='  -6.97%  ' }
    @{ Row=28; D='20.29'; E='  +2.48%  ' }
    @{ Row=29; D='145.34'; E='  -0.32%  ' }
    @{ Row=30; D='5.070'; E='  +2.85%  ' }
    @{ Row=31; D='126.18'; E='  +0.90%  ' }
    @{ Row=32; D='1.748.92'; E='  +0.29%  ' }
    @{ Row=33; D='6.206'; E='  -0.91%  ' }
    @{ Row=34; E='  +3.72%  ' }
    @{ Row=35; D='2.021'; E='  +2.08%  ' }
    @{ Row=36; D='10.12'; E='  -2.72%  ' }
    @{ Row=37; D='0.08565'; E='  +1.38%  ' }
    @{ Row=38; D='0.02560'; E='  +1.03%  ' }
    @{ Row=39; E='  +0.88%  ' }
    @{ Row=40; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.06539'; E='  +0.27%  ' }
    @{ Row=41; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='1.330'; E='  +5.60%  ' }
    @{ Row=42; D='5.476'; E='  -0.34%  ' }
    @{ Row=43; D='11.63'; E='  -1.03%  ' }
    @{ Row=44; D='0.6469'; E='  +1.54%  ' }
    @{ Row=45; D='14.09'; E='  -2.29%  ' }
    @{ Row=46; E='  +0.14%  ' }
    @{ Row=47; D='0.6035'; E='  +0.78%  ' }
    @{ Row=48; D='3.793'; E='  +0.42%  ' }
    @{ Row=49; D='1.306'; E='  +10.27%  ' }
    @{ Row=50; D='2.096'; E='  -1.02%  ' }
    @{ Row=51; D='125.83'; E='  +3.12%  ' }
)

foreach ($u in $updates) {
    foreach ($col in @('B','C','D','E')) {
        if ($u.ContainsKey($col)) {
            $cellRef = "$col$($u.Row)"
            $newVal = $u[$col]
            $range = $ws.Range($cellRef)
            if ($newVal -match '^[+-]?\d+(\.\d+)?$') {
                # Value looks like a plain number (e.g. "49.90" or "0.3415").
                # Force Text storage first so Excel keeps the exact original
                # string (trailing zeros / fixed decimals) instead of silently
                # coercing it into a floating point number, then clear the
                # temporary format so no stray style is left on the cell.
                $range.NumberFormat = "@"
                $range.Value = $newVal
                $range.ClearFormats()
            } else {
                $range.Value = $newVal
            }
        }
    }
}
